$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.388571333333333
$ws.Range("H2").Value = 4.165713999999999
$ws.Range("I2").Value = 0.3523526610542377
$ws.Range("J2").Value = 0.3523526610542378
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 58.88086666666666
$ws.Range("N2").Value = 176.6426
$ws.Range("O2").Value = 0.2818386429293818
$ws.Range("P2").Value = 0.2818386429293819
$ws.Range("Q2").Value = 81.76028353515554
$ws.Range("R2").Value = 735.8425518163998
$ws.Range("S2").Value = 0.09930659582408281
$ws.Range("T2").Value = 0.09930659582408284

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.388571333333333
$ws.Range("H3").Value = 4.165713999999999
$ws.Range("I3").Value = 0.3523526610542377
$ws.Range("J3").Value = 0.3523526610542378
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 110.7417803333333
$ws.Range("N3").Value = 332.225341
$ws.Range("O3").Value = 0.5300756400448766
$ws.Range("P3").Value = 0.5300756400448767
$ws.Range("Q3").Value = 153.7728615731638
$ws.Range("R3").Value = 1383.955754158474
$ws.Range("S3").Value = 0.1867735623298405
$ws.Range("T3").Value = 0.1867735623298406

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.388571333333333
$ws.Range("H4").Value = 4.165713999999999
$ws.Range("I4").Value = 0.3523526610542377
$ws.Range("J4").Value = 0.3523526610542378
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 39.29429233333334
$ws.Range("N4").Value = 117.882877
$ws.Range("O4").Value = 0.1880857170257415
$ws.Range("P4").Value = 0.1880857170257415
$ws.Range("Q4").Value = 54.56292789768644
$ws.Range("R4").Value = 491.066351079178
$ws.Range("S4").Value = 0.06627250290031436
$ws.Range("T4").Value = 0.06627250290031438

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.552285333333333
$ws.Range("H5").Value = 7.656856
$ws.Range("I5").Value = 0.6476473389457622
$ws.Range("J5").Value = 0.6476473389457623
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 58.88086666666666
$ws.Range("N5").Value = 176.6426
$ws.Range("O5").Value = 0.2818386429293818
$ws.Range("P5").Value = 0.2818386429293819
$ws.Range("Q5").Value = 150.2807724072889
$ws.Range("R5").Value = 1352.5269516656
$ws.Range("S5").Value = 0.182532047105299
$ws.Range("T5").Value = 0.1825320471052991

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.552285333333333
$ws.Range("H6").Value = 7.656856
$ws.Range("I6").Value = 0.6476473389457622
$ws.Range("J6").Value = 0.6476473389457623
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 110.7417803333333
$ws.Range("N6").Value = 332.225341
$ws.Range("O6").Value = 0.5300756400448766
$ws.Range("P6").Value = 0.5300756400448767
$ws.Range("Q6").Value = 282.6446217319885
$ws.Range("R6").Value = 2543.801595587896
$ws.Range("S6").Value = 0.3433020777150361
$ws.Range("T6").Value = 0.3433020777150362

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.552285333333333
$ws.Range("H7").Value = 7.656856
$ws.Range("I7").Value = 0.6476473389457622
$ws.Range("J7").Value = 0.6476473389457623
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 39.29429233333334
$ws.Range("N7").Value = 117.882877
$ws.Range("O7").Value = 0.1880857170257415
$ws.Range("P7").Value = 0.1880857170257415
$ws.Range("Q7").Value = 100.2902460060791
$ws.Range("R7").Value = 902.6122140547121
$ws.Range("S7").Value = 0.1218132141254271
$ws.Range("T7").Value = 0.1218132141254271

